$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 5
$ws.Range("D2").Value = -47.200000000000003
$ws.Range("F2").Value = -47.200000000000003
$ws.Range("H2").Value = 0.49782812780356078
$ws.Range("I2").Value = 0.4062696080437419
$ws.Range("J2").Value = 0.48665956450759856

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("D3").Value = -47.200000000000003
$ws.Range("F3").Value = -47.200000000000003
$ws.Range("H3").Value = 1.0936890868282563
$ws.Range("I3").Value = 0.52968928668477566
$ws.Range("J3").Value = 0.39669952206598796

# Row 4
$ws.Range("A4").Value = 7
$ws.Range("D4").Value = -47.200000000000003
$ws.Range("F4").Value = -47.200000000000003
$ws.Range("H4").Value = 1.0936890868282563
$ws.Range("I4").Value = 0.52968928668477566
$ws.Range("J4").Value = 0.39669952206598796

# Row 5
$ws.Range("A5").Value = 8
$ws.Range("D5").Value = -48.800000000000004
$ws.Range("F5").Value = -48.800000000000004
$ws.Range("H5").Value = 1.181359437926748
$ws.Range("I5").Value = 0.067205378687909836
$ws.Range("J5").Value = 0.68802191628728926

# Row 6
$ws.Range("A6").Value = 2
$ws.Range("H6").Value = 1.181359437926748
$ws.Range("I6").Value = 0.067205378687909836
$ws.Range("J6").Value = 0.68802191628728926

# Row 7
$ws.Range("A7").Value = 4
$ws.Range("D7").Value = -49.600000000000001
$ws.Range("H7").Value = 0.81288935844146215
$ws.Range("I7").Value = 0.48566544174909315
$ws.Range("J7").Value = 0.6874499523968628

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("D8").Value = -40
$ws.Range("H8").Value = 0.35282568648548157
$ws.Range("I8").Value = 0.50813390051494534
$ws.Range("J8").Value = 0.073974240930933854

# Row 9
$ws.Range("A9").Value = 3
$ws.Range("D9").Value = -48.800000000000004
$ws.Range("H9").Value = 0.99767385091727534
$ws.Range("I9").Value = 0.38954505945686929
$ws.Range("J9").Value = 0.6338666319655698

# Row 10
$ws.Range("D10").Value = -47.200000000000003
$ws.Range("E10").Value = 1550.1412127916667
$ws.Range("F10").Value = -47.200000000000003
$ws.Range("H10").Value = 1.0936890868282563
$ws.Range("I10").Value = 0.52968928668477566
$ws.Range("J10").Value = 0.39669952206598796

# Row 11
$ws.Range("E11").Value = 84.321910958333334
$ws.Range("F11").Value = -48.800000000000004
$ws.Range("H11").Value = 1.181359437926748
$ws.Range("I11").Value = 0.067205378687909836
$ws.Range("J11").Value = 0.68802191628728926

# Row 12
$ws.Range("D12").Value = -48.800000000000004
$ws.Range("E12").Value = 12.172907541666667
$ws.Range("F12").Value = -48.800000000000004
$ws.Range("H12").Value = 0.99767385091727534
$ws.Range("I12").Value = 0.38954505945686929
$ws.Range("J12").Value = 0.6338666319655698

# Row 13
$ws.Range("D13").Value = -49.600000000000001
$ws.Range("E13").Value = 90.518051958333331
$ws.Range("H13").Value = 0.81288935844146215
$ws.Range("I13").Value = 0.48566544174909315
$ws.Range("J13").Value = 0.6874499523968628

# Row 14
$ws.Range("D14").Value = -47.200000000000003
$ws.Range("E14").Value = 185.3418795
$ws.Range("H14").Value = 0.49782812780356078
$ws.Range("I14").Value = 0.4062696080437419
$ws.Range("J14").Value = 0.48665956450759856

# Column width adjustments (H and J narrow from 14.7109375 to 13.7109375)
$ws.Columns.Item(8).ColumnWidth = 12.8333325
$ws.Columns.Item(10).ColumnWidth = 12.8333325
